$d = $word.ActiveDocument

# 1) "{% foreach field in fields.Group %} " -> "{% for field in fields.Group %} "
$d.Content.Find.Execute("foreach field in fields.Group", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "for field in fields.Group", 2)

# 2) "{% endeach %}" -> "{% end" | "for" | <bookmark _GoBack> | "%}" | " "
#    Locate the paragraph that still holds the "{% endeach %}" marker.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("{% endeach %}")) {
        $target = $para
    }
}
$paraStart = $target.Range.Start

# "{% endeach %} "
#  0123456789...
# Replace "each " (the word plus the following space) with "for" (no trailing space),
# turning "{% endeach %} " into "{% endfor%} " (still a single run at this point).
$eachRange = $d.Range($paraStart + 6, $paraStart + 11)
$eachRange.Text = "for"

# Split "{% end" / "for" apart by binding+removing a throw-away bookmark at the boundary
# (removing a bookmark does not re-merge runs that were split to host it).
$b1 = $d.Range($paraStart + 6, $paraStart + 6)
$d.Bookmarks.Add("zSplitTmp1", $b1)
$d.Bookmarks("zSplitTmp1").Delete()

# Drop the real "_GoBack" bookmark exactly between "for" and "%}" (it used to live in the
# next-but-one empty paragraph; moving it here removes it from there automatically since a
# bookmark name is unique document-wide).
$b2 = $d.Range($paraStart + 9, $paraStart + 9)
$d.Bookmarks.Add("_GoBack", $b2)

# Split "%}" / " " apart the same throw-away-bookmark way
$b3 = $d.Range($paraStart + 11, $paraStart + 11)
$d.Bookmarks.Add("zSplitTmp2", $b3)
$d.Bookmarks("zSplitTmp2").Delete()
